# Swap the values of the two-row "~UC_Sets:" header blocks (A-column)
# back to their originally-intended order on each relevant worksheet.
# Pattern: for each pair of rows (n, n+1), the value that was in A(n)
# moves to A(n+1), and the value that was in A(n+1) moves to A(n).

$wb = $excel.ActiveWorkbook

function Swap-ARows {
    param($ws, [int]$row1, [int]$row2)

    $cellA = $ws.Cells.Item($row1, 1)
    $cellB = $ws.Cells.Item($row2, 1)

    $valA = $cellA.Value()
    $valB = $cellB.Value()

    $cellA.Value = $valB
    $cellB.Value = $valA
}

# Sheet "Cars": rows 1-2 and 7-8
$ws = $wb.Worksheets.Item("Cars")
Swap-ARows $ws 1 2
Swap-ARows $ws 7 8

# Sheet "Cars_2020": rows 1-2
$ws = $wb.Worksheets.Item("Cars_2020")
Swap-ARows $ws 1 2

# Sheet "CCS+h2": rows 1-2
$ws = $wb.Worksheets.Item("CCS+h2")
Swap-ARows $ws 1 2

# Sheet "CH_RH": rows 1-2
$ws = $wb.Worksheets.Item("CH_RH")
Swap-ARows $ws 1 2

# Sheet "IND_fuels": rows 1-2
$ws = $wb.Worksheets.Item("IND_fuels")
Swap-ARows $ws 1 2

# Sheet "Power_sector": rows 1-2, 10-11, 17-18
$ws = $wb.Worksheets.Item("Power_sector")
Swap-ARows $ws 1 2
Swap-ARows $ws 10 11
Swap-ARows $ws 17 18

# Sheet "Thermal_gencap": rows 1-2 and 7-8
$ws = $wb.Worksheets.Item("Thermal_gencap")
Swap-ARows $ws 1 2
Swap-ARows $ws 7 8

# Sheet "TRA_Policy": rows 1-2
$ws = $wb.Worksheets.Item("TRA_Policy")
Swap-ARows $ws 1 2
